$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, pushing the existing rows 70-91 down to 71-92.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly price record.
$ws.Cells.Item(70, 1).Value  = 10
$ws.Cells.Item(70, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value  = "La Araucanía"
$ws.Cells.Item(70, 4).Value  = 44782
$ws.Cells.Item(70, 5).Value  = 9
$ws.Cells.Item(70, 6).Value  = "Fruta"
$ws.Cells.Item(70, 7).Value  = 100108
$ws.Cells.Item(70, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(70, 9).Value  = 100108007
$ws.Cells.Item(70, 10).Value = "Coco"
$ws.Cells.Item(70, 11).Value = "Sin especificar"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 20
$ws.Cells.Item(70, 14).Value = 30000
$ws.Cells.Item(70, 15).Value = 30000
$ws.Cells.Item(70, 16).Value = 30000
$ws.Cells.Item(70, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(70, 18).Value = "Perú"
$ws.Cells.Item(70, 19).Value = 1500
$ws.Cells.Item(70, 20).Value = 20
